# cryptos list refresh — update Price (D) and Volume(1h) (E) columns
# with latest scraped values; two rows (46/47) also swap identity
# (Kaspa <-> FraxShare) because the upstream ranking API reordered them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.353.27'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.805.42'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.575'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.96%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.35'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.61%  '
$ws.Range('E9').Value = '  +2.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0695'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0962'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.98%  '
$ws.Range('D12').Value = '2.065.33'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.88'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.80%  '
$ws.Range('D14').Value = '1.827.75'
$ws.Range('E14').Value = '  +3.31%  '
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.53'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.09%  '
$ws.Range('D17').Value = '34.341.86'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').Value = '0.0₃0797'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.22%  '
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '171.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.13'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.04%  '
$ws.Range('E26').Value = '  +9.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.119'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0533'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('E33').Value = '  +1.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('D35').Value = '1.397.67'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.675'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.08%  '
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0192'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.968'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.68%  '
$ws.Range('E41').Value = '  +10.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '82.83'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.43'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0506'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.60%  '
$ws.Range('D48').Value = '1.966.37'
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.68'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('E51').Value = '  +0.20%  '
